# Append a new scraped-job row (2025-09-17 18:25:35 JST run) to the
# "ランサーズ" sheet, and refresh the "取得日時" timestamp on every
# existing row to match the new scrape run.
#
# Source data gained one new posting that now sorts into row 7 by
# priority score (125), pushing everything that used to be row 7..22
# down to row 8..23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-17 18:25:35"

# 1) Insert a fresh row at position 7 - this shifts the existing row 7..22
#    (and their cell formatting) down to row 8..23, and grows the used
#    range/dimension automatically.
$ws.Rows(7).Insert()

# 2) Refresh the timestamp in column A for every data row (2..23) so the
#    whole sheet reflects the new scrape run.
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# 3) Populate the brand-new row 7 with the newly scraped job posting.
$ws.Range("B7").Value = "初回 Laravel Livewireを使ったWebシステム開発の募集"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5395502"
$ws.Range("G7").Value = 125
$ws.Range("H7").Value = "◆開発,システム開発"

# 4) The engine's Rows.Insert() does not shift the sheet's <hyperlinks>
#    table, so rebuild it from scratch: wipe whatever is left (deleting
#    any single hyperlink clears the whole collection here) and re-add
#    one hyperlink per URL cell in column F, in row order, using each
#    cell's own (already-correct) text as the link target.
$ws.Range("F2").Hyperlinks.Delete()

for ($r = 2; $r -le 23; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value2
    $ws.Hyperlinks.Add($cell, $url)
    # Hyperlinks.Add() re-applies the "Hyperlink" cell style but as a fresh
    # duplicate style record; re-assert the named style explicitly so it
    # resolves back to the workbook's existing single "Hyperlink" style
    # entry instead of growing a near-duplicate one.
    $cell.Style = "Hyperlink"
}
